# Applies the scheduled-runner profit/price refresh to the Leviathan_Profits
# workbook: updates currentAveragePrice / LevePrice / LeveProfit columns
# (H, I, J, K, L, M, N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# with freshly scraped market-board figures. Some rows gain newly-populated
# profit cells (previously blank) while a couple of rows lose a stale
# profit cell that's superseded by the refreshed value.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 4591.4614
$ws.Range("J29").Value = 6488.1113
$ws.Range("L29").Value = 19464.3339
$ws.Range("N29").Value = -20026.3339
$ws.Range("H62").Value = 11999
$ws.Range("J62").Value = 11999
$ws.Range("L62").Value = 11999
$ws.Range("N62").Value = -13247
$ws.Range("H65").Value = 11999
$ws.Range("J65").Value = 11999
$ws.Range("L65").Value = 59995
$ws.Range("N65").Value = -66235
$ws.Range("H70").Value = 3678.0557
$ws.Range("I70").Value = 3338.4285
$ws.Range("J70").Value = 3894.182
$ws.Range("K70").Value = 10015.2855
$ws.Range("L70").Value = 11682.546
$ws.Range("M70").Value = -9745.2855
$ws.Range("N70").Value = -12222.546
$ws.Range("H73").Value = 3678.0557
$ws.Range("I73").Value = 3338.4285
$ws.Range("J73").Value = 3894.182
$ws.Range("K73").Value = 10015.2855
$ws.Range("L73").Value = 11682.546
$ws.Range("M73").Value = -9079.2855
$ws.Range("N73").Value = -13554.546
$ws.Range("H113").Value = 81776.38
$ws.Range("I113").Value = 146299
$ws.Range("K113").Value = 146299
$ws.Range("M113").Value = -143045
$ws.Range("H116").Value = 4979.091
$ws.Range("I116").Value = 4945
$ws.Range("K116").Value = 4945
$ws.Range("M116").Value = -1503
$ws.Range("H125").Value = 1313.4
$ws.Range("I125").Value = 878
$ws.Range("J125").Value = 1966.5
$ws.Range("K125").Value = 7902
$ws.Range("L125").Value = 17698.5
$ws.Range("M125").Value = -5442
$ws.Range("N125").Value = -22618.5
$ws.Range("H127").Value = 328
$ws.Range("I127").Value = 328
$ws.Range("K127").Value = 984
$ws.Range("M127").Value = 3976
$ws.Range("H138").Value = 1729.9579
$ws.Range("J138").Value = 1879.919
$ws.Range("L138").Value = 5639.757000000001
$ws.Range("N138").Value = -15919.757
$ws.Range("H141").Value = 3189.3
$ws.Range("I141").Value = 3210.3333
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 9630.999899999999
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = -4450.999899999999
$ws.Range("N141").Value = -19360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2232.5386
$ws.Range("I61").Value = 1187.875
$ws.Range("K61").Value = 1187.875
$ws.Range("M61").Value = -975.875
$ws.Range("H110").Value = 1543.5
$ws.Range("J110").Value = 1413
$ws.Range("L110").Value = 1413
$ws.Range("N110").Value = -5503
$ws.Range("H136").Value = 2232.5386
$ws.Range("I136").Value = 1187.875
$ws.Range("K136").Value = 3563.625
$ws.Range("M136").Value = -1013.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3016.6316
$ws.Range("I86").Value = 3487.8
$ws.Range("K86").Value = 3487.8
$ws.Range("M86").Value = -2364.8
$ws.Range("H89").Value = 3016.6316
$ws.Range("I89").Value = 3487.8
$ws.Range("K89").Value = 17439
$ws.Range("M89").Value = -11823
$ws.Range("H134").Value = 1539.9688
$ws.Range("I134").Value = 1416.2142
$ws.Range("K134").Value = 4248.642599999999
$ws.Range("M134").Value = -1713.642599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 777.55554
$ws.Range("I22").Value = 428.2857
$ws.Range("K22").Value = 428.2857
$ws.Range("M22").Value = -78.28570000000002
$ws.Range("H50").Value = 14999.75
$ws.Range("J50").Value = 14999.75
$ws.Range("L50").Value = 14999.75
$ws.Range("N50").Value = -16249.75
$ws.Range("H51").Value = 14999.857
$ws.Range("J51").Value = 14999.857
$ws.Range("L51").Value = 14999.857
$ws.Range("N51").Value = -16471.857
$ws.Range("H59").Value = 19998.428
$ws.Range("I59").Value = 19990
$ws.Range("J59").Value = 19999.834
$ws.Range("K59").Value = 19990
$ws.Range("L59").Value = 19999.834
$ws.Range("M59").Value = -18845
$ws.Range("N59").Value = -22289.834
$ws.Range("H60").Value = 14999.875
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 14999.875
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 14999.875
$ws.Range("M60").Value = ""
$ws.Range("N60").Value = -16021.875
$ws.Range("H61").Value = 14999.857
$ws.Range("J61").Value = 14999.857
$ws.Range("L61").Value = 14999.857
$ws.Range("N61").Value = -15695.857
$ws.Range("H62").Value = 4800
$ws.Range("J62").Value = 4800
$ws.Range("L62").Value = 4800
$ws.Range("N62").Value = -6048
$ws.Range("H65").Value = 4800
$ws.Range("J65").Value = 4800
$ws.Range("L65").Value = 24000
$ws.Range("N65").Value = -30240
$ws.Range("H99").Value = 18429.777
$ws.Range("I99").Value = 46318
$ws.Range("J99").Value = 4485.6665
$ws.Range("K99").Value = 46318
$ws.Range("L99").Value = 4485.6665
$ws.Range("M99").Value = -44820
$ws.Range("N99").Value = -7481.6665
$ws.Range("H126").Value = 18429.777
$ws.Range("I126").Value = 46318
$ws.Range("J126").Value = 4485.6665
$ws.Range("K126").Value = 138954
$ws.Range("L126").Value = 13456.9995
$ws.Range("M126").Value = -136484
$ws.Range("N126").Value = -18396.9995
$ws.Range("H132").Value = 3295.44
$ws.Range("I132").Value = 3762.05
$ws.Range("K132").Value = 11286.15
$ws.Range("M132").Value = -8756.150000000001
$ws.Range("H134").Value = 1824.5938
$ws.Range("I134").Value = 1271.9656
$ws.Range("J134").Value = 7166.6665
$ws.Range("K134").Value = 3815.8968
$ws.Range("L134").Value = 21499.9995
$ws.Range("M134").Value = -1280.8968
$ws.Range("N134").Value = -26569.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 163
$ws.Range("J12").Value = 194.5
$ws.Range("L12").Value = 583.5
$ws.Range("N12").Value = -929.5
$ws.Range("H22").Value = 1500
$ws.Range("I22").Value = 1500
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 4500
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -4331
$ws.Range("N22").Value = ""
$ws.Range("H27").Value = 1500
$ws.Range("I27").Value = 1500
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 4500
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -4398
$ws.Range("N27").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 35716456
$ws.Range("I107").Value = 344.625
$ws.Range("J107").Value = 83337940
$ws.Range("K107").Value = 344.625
$ws.Range("L107").Value = 83337940
$ws.Range("M107").Value = 1575.375
$ws.Range("N107").Value = -83341780
$ws.Range("H113").Value = 3711.476
$ws.Range("I113").Value = 2747.7
$ws.Range("J113").Value = 4587.636
$ws.Range("K113").Value = 2747.7
$ws.Range("L113").Value = 4587.636
$ws.Range("M113").Value = -577.6999999999998
$ws.Range("N113").Value = -8927.636
$ws.Range("H132").Value = 1985.475
$ws.Range("I132").Value = 2022.1111
$ws.Range("K132").Value = 6066.3333
$ws.Range("M132").Value = -3536.3333
$ws.Range("H136").Value = 24042.227
$ws.Range("J136").Value = 24042.227
$ws.Range("L136").Value = 72126.681
$ws.Range("N136").Value = -77226.681

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 12712.444
$ws.Range("I7").Value = 16185.842
$ws.Range("K7").Value = 16185.842
$ws.Range("M7").Value = -16073.842
$ws.Range("H61").Value = 60839.75
$ws.Range("I61").Value = 46961.5
$ws.Range("K61").Value = 46961.5
$ws.Range("M61").Value = -46759.5
$ws.Range("H113").Value = 60839.75
$ws.Range("I113").Value = 46961.5
$ws.Range("K113").Value = 46961.5
$ws.Range("M113").Value = -44791.5
$ws.Range("H126").Value = 12712.444
$ws.Range("I126").Value = 16185.842
$ws.Range("K126").Value = 48557.526
$ws.Range("M126").Value = -46087.526
$ws.Range("H129").Value = 85281.336
$ws.Range("J129").Value = 85281.336
$ws.Range("L129").Value = 85281.336
$ws.Range("N129").Value = -95281.336
$ws.Range("H130").Value = 91500
$ws.Range("I130").Value = 93000
$ws.Range("J130").Value = 90000
$ws.Range("K130").Value = 93000
$ws.Range("L130").Value = 90000
$ws.Range("M130").Value = -87980
$ws.Range("N130").Value = -100040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 23975
$ws.Range("I62").Value = 7950
$ws.Range("K62").Value = 7950
$ws.Range("M62").Value = -7326
$ws.Range("H65").Value = 23975
$ws.Range("I65").Value = 7950
$ws.Range("K65").Value = 39750
$ws.Range("M65").Value = -36630
$ws.Range("H70").Value = 18143.691
$ws.Range("J70").Value = 17551.916
$ws.Range("L70").Value = 17551.916
$ws.Range("N70").Value = -18181.916
$ws.Range("H73").Value = 18143.691
$ws.Range("J73").Value = 17551.916
$ws.Range("L73").Value = 17551.916
$ws.Range("N73").Value = -19735.916
$ws.Range("H105").Value = 35307.5
$ws.Range("J105").Value = 35307.5
$ws.Range("L105").Value = 35307.5
$ws.Range("N105").Value = -42295.5
$ws.Range("H122").Value = 2128.9
$ws.Range("I122").Value = 2042.375
$ws.Range("K122").Value = 6127.125
$ws.Range("M122").Value = -3677.125
$ws.Range("H128").Value = 59857.5
$ws.Range("J128").Value = 59857.5
$ws.Range("L128").Value = 59857.5
$ws.Range("N128").Value = -69817.5
